$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record at row 75 ("Cucha Cucha 2918" / Aplomar), pushing
# every existing row from 75 downward down by one (old row 92 becomes row 93).
$ws.Range("A75").EntireRow.Insert()

# Columns that look numeric/date-like (Caso, F. De Reclamo, Comuna, OT) must
# stay plain text - matches every other data row in the sheet, which stores
# these as literal strings even when they look numeric/date-like.
$ws.Range("A75:B75").NumberFormat = "@"
$ws.Range("D75:E75").NumberFormat = "@"

$ws.Range("A75").Value = "807044223"
$ws.Range("B75").Value = "5/22/2025"
$ws.Range("C75").Value = "Cucha Cucha 2918"
$ws.Range("D75").Value = "7"
$ws.Range("E75").Value = "807044223"
$ws.Range("F75").Value = "AYKO"
$ws.Range("G75").Value = "Pendiente"
$ws.Range("H75").Value = "Aplomar"
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = "Aplomo"
$ws.Range("K75").Value = "Sin equipos"
$ws.Range("L75").Value = "Pasante"
$ws.Range("M75").Value = -58.469783
$ws.Range("N75").Value = -34.599214
$ws.Range("O75").Value = "Paternal"
$ws.Range("P75").Value = "Capital Norte"
